$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: B1 "Court order " (trailing space) -> "Court order", styled like the
# --- Client Name data cells (Arial font, matching style used in column C) ---
$ws.Range("B1").Value = "Court order"
$ws.Range("B1").Font.Name = "Arial"

# --- Row 5: update case number, court order number and client name ---
$ws.Range("A5").Value = 11052066
$ws.Range("B5").Value = 1
$ws.Range("C5").Value = "CHARLIE THOMPSON"

# --- New row 6: add another case for LEE OSWARLD ---
$ws.Range("A6").Value = 11052077
$ws.Range("B6").Value = 2
$ws.Range("C6").Value = "LEE OSWARLD "
$ws.Range("C6").Font.Name = "Arial"

$ws.Range("D6").Value = 45521
$ws.Range("E6").Value = 45579

# Match date formatting of the row above for the new date cells
$ws.Range("D5").Copy($ws.Range("D6"))
$ws.Range("E5").Copy($ws.Range("E6"))
$ws.Range("D6").Value = 45521
$ws.Range("E6").Value = 45579

# Match row height formatting used by the rest of the data rows
$ws.Rows(6).RowHeight = 15.75

# Update the active selection to reflect where the user ended up after editing
$ws.Range("A6").Select() | Out-Null
